$wb = $excel.ActiveWorkbook

# --- Duplicate "BasicAddStat" to create the new "BasicMulStat" sheet,
#     placed right after BasicAddStat (i.e. before Storage). This also
#     carries over the header row, styles, and the B1 cell comment. ---
$srcSheet = $wb.Worksheets.Item("BasicAddStat")
$srcSheet.Copy($null, $srcSheet) | Out-Null
$mulSheet = $wb.Worksheets.Item("BasicAddStat (2)")
$mulSheet.Name = "BasicMulStat"

# --- Rework the copied data into the "multiplier" variant of the table ---
$mulSheet.Range("B2:B5").Value = 1.1
$mulSheet.Range("E2").Value = 0.1
$mulSheet.Range("G2:H3").Value = 0.1
$mulSheet.Range("I2").Value = 0.1
$mulSheet.Range("J2").Value = 0.1
$mulSheet.Range("K2").Value = 0.1
$mulSheet.Range("C4").Value = 0.1
$mulSheet.Range("D5").Value = 1

# --- Match the printed-page setup carried over onto the new sheet ---
$mulSheet.PageSetup.PaperSize = 9
$mulSheet.PageSetup.Orientation = 1

# --- Restore per-sheet selections as saved in the workbook ---
$itemSheet = $wb.Worksheets.Item("Item")
$itemSheet.Range("D46").Select() | Out-Null

$addSheet = $wb.Worksheets.Item("BasicAddStat")
$addSheet.Range("F28").Select() | Out-Null

$mulSheet.Range("C4").Select() | Out-Null

$storageSheet = $wb.Worksheets.Item("Storage")
$storageSheet.Range("D42").Select() | Out-Null

# --- "BasicAddStat" is the active/visible tab when the workbook is saved ---
$addSheet.Activate() | Out-Null
$addSheet.Range("F28").Select() | Out-Null
